$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "longitude" and "latitude" header labels in D1 and E1
$d1 = $ws.Range("D1").Value()
$e1 = $ws.Range("E1").Value()
$ws.Range("D1").Value = $e1
$ws.Range("E1").Value = $d1

# Move the active selection to E1 (reflects final cursor position in the diff)
$ws.Range("E1").Select()

# Header row now wraps to two lines with the new label, matching the
# taller row height Excel settles on after the edit
$ws.Rows.Item(1).RowHeight = 29

